$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 4 new rows before row 12 (the "feature set - balanced data set"
#    header), to make room for the new v3 / v4 cognatehood results that are
#    appended under the "unbalanced" table (rows 10-13).
#    This pushes the whole "balanced" table (old rows 12-20) down to rows 16-24.
# ---------------------------------------------------------------------------
$ws.Rows("10:13").Insert()

# ---------------------------------------------------------------------------
# 2. Highlight (bold) the best unbalanced-set result, row 7 (v1 + words(100))
# ---------------------------------------------------------------------------
$ws.Range("B7").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. Fill in the new "unbalanced" rows 10-13: cognatehood v3 / v4 results
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "pos tag + len + cognatehood (v3)"
$ws.Range("B10").Value = 0.86
$ws.Range("C10").Formula = "=((1-B2)-(1-B10)) / (1-B2) * 100"

$ws.Range("A11").Value = "pos tag + len + cognatehood (v3) + words (100)"
$ws.Range("B11").Value = 0.86399999999999999
$ws.Range("B11").Font.Bold = $true
$ws.Range("C11").Formula = "=((1-B2)-(1-B11)) / (1-B2) * 100"

$ws.Range("A12").Value = "pos tag + len + cognatehood (v4)"
$ws.Range("B12").Value = 0.85899999999999999
$ws.Range("C12").Formula = "=((1-B2)-(1-B12)) / (1-B2) * 100"

$ws.Range("A13").Value = "pos tag + len + cognatehood (v4) + words (100)"
$ws.Range("B13").Value = 0.86299999999999999
$ws.Range("C13").Formula = "=((1-B2)-(1-B13)) / (1-B2) * 100"

# ---------------------------------------------------------------------------
# 4. Highlight (bold) the best balanced-set result, now at row 22
#    (old row 18, v1 + words(100))
# ---------------------------------------------------------------------------
$ws.Range("B22").Font.Bold = $true

# ---------------------------------------------------------------------------
# 5. Append the new "balanced" rows 25-28: cognatehood v3 / v4 results
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "pos tag + len + cognatehood (v3)"
$ws.Range("B25").Value = 0.74199999999999999
$ws.Range("C25").Formula = "=((1-B17)-(1-B25)) / (1-B17) * 100"

$ws.Range("A26").Value = "pos tag + len + cognatehood (v3) + words (100)"
$ws.Range("B26").Value = 0.747
$ws.Range("C26").Formula = "=((1-B17)-(1-B26)) / (1-B17) * 100"

$ws.Range("A27").Value = "pos tag + len + cognatehood (v4)"
$ws.Range("B27").Value = 0.73799999999999999
$ws.Range("C27").Formula = "=((1-B17)-(1-B27)) / (1-B17) * 100"

$ws.Range("A28").Value = "pos tag + len + cognatehood (v4) + words (100)"
$ws.Range("B28").Value = 0.745
$ws.Range("C28").Formula = "=((1-B17)-(1-B28)) / (1-B17) * 100"

# ---------------------------------------------------------------------------
# 6. Update the two charts that plot the "balanced" table, whose source rows
#    moved from 12-20 down to 16-24.
# ---------------------------------------------------------------------------
$balancedChart = $ws.ChartObjects().Item(2).Chart
$balancedChart.SeriesCollection().Item(1).Formula = "=SERIES(Sheet1!`$B`$16,Sheet1!`$A`$17:`$A`$24,Sheet1!`$B`$17:`$B`$24,1)"
$balancedChart.SeriesCollection().Item(2).Formula = "=SERIES(Sheet1!`$C`$16,Sheet1!`$A`$17:`$A`$24,Sheet1!`$C`$17:`$C`$24,2)"

# ---------------------------------------------------------------------------
# 7. Track the chart drawing anchors with the newly inserted rows (default
#    row height = 14.4pt, 4 rows were inserted).
#    - Chart 1 (the "un balanced" chart) spans across the insertion point
#      (rows 0-18), so it grows taller (bottom edge pushed down) while its
#      top stays put.
#    - Chart 2 (the "Balanced" chart) sits entirely below the insertion
#      point (rows 29-51), so the whole chart simply shifts down.
# ---------------------------------------------------------------------------
$rowShiftPts = 4 * 14.4
$chartObj1 = $ws.ChartObjects().Item(1)
$chartObj1.Height = $chartObj1.Height + $rowShiftPts
$chartObj2 = $ws.ChartObjects().Item(2)
$chartObj2.Top = $chartObj2.Top + $rowShiftPts

# ---------------------------------------------------------------------------
# 8. View state: zoom to 115% and select the full used range, ending on C28.
# ---------------------------------------------------------------------------
$ws.Select()
$excel.ActiveWindow.Zoom = 115
$ws.Range("A1:C28").Select()
